$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 3907.625
$ws.Range("I31").Value = 4351.5713
$ws.Range("J31").Value = 800
$ws.Range("K31").Value = 13054.7139
$ws.Range("L31").Value = 2400
$ws.Range("M31").Value = -12824.7139
$ws.Range("N31").Value = -2860
$ws.Range("H51").Value = 2980.1094
$ws.Range("J51").Value = 2952.9412
$ws.Range("L51").Value = 2952.9412
$ws.Range("N51").Value = -3920.9412
$ws.Range("H76").Value = 5234.3076
$ws.Range("I76").Value = 4002.8
$ws.Range("K76").Value = 4002.8
$ws.Range("M76").Value = -3687.8
$ws.Range("H79").Value = 5234.3076
$ws.Range("I79").Value = 4002.8
$ws.Range("K79").Value = 4002.8
$ws.Range("M79").Value = -2910.8
$ws.Range("H128").Value = 150000
$ws.Range("J128").Value = 150000
$ws.Range("L128").Value = 150000
$ws.Range("N128").Value = -159960
$ws.Range("H132").Value = 15057.786
$ws.Range("I132").Value = 22780.223
$ws.Range("J132").Value = 1157.4
$ws.Range("K132").Value = 68340.66900000001
$ws.Range("L132").Value = 3472.2
$ws.Range("M132").Value = -65810.66900000001
$ws.Range("N132").Value = -8532.200000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3462.5356
$ws.Range("I61").Value = 3152.261
$ws.Range("J61").Value = 4889.8
$ws.Range("K61").Value = 3152.261
$ws.Range("L61").Value = 4889.8
$ws.Range("M61").Value = -2940.261
$ws.Range("N61").Value = -5313.8
$ws.Range("H136").Value = 3462.5356
$ws.Range("I136").Value = 3152.261
$ws.Range("J136").Value = 4889.8
$ws.Range("K136").Value = 9456.782999999999
$ws.Range("L136").Value = 14669.4
$ws.Range("M136").Value = -6906.782999999999
$ws.Range("N136").Value = -19769.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 37265.83
$ws.Range("I20").Value = 62024.293
$ws.Range("J20").Value = 2191.3333
$ws.Range("K20").Value = 62024.293
$ws.Range("L20").Value = 2191.3333
$ws.Range("M20").Value = -61777.293
$ws.Range("N20").Value = -2685.3333
$ws.Range("H86").Value = 5499.1665
$ws.Range("I86").Value = 6123.5
$ws.Range("J86").Value = 4250.5
$ws.Range("K86").Value = 6123.5
$ws.Range("L86").Value = 4250.5
$ws.Range("M86").Value = -5000.5
$ws.Range("N86").Value = -6496.5
$ws.Range("H89").Value = 5499.1665
$ws.Range("I89").Value = 6123.5
$ws.Range("J89").Value = 4250.5
$ws.Range("K89").Value = 30617.5
$ws.Range("L89").Value = 21252.5
$ws.Range("M89").Value = -25001.5
$ws.Range("N89").Value = -32484.5
$ws.Range("H107").Value = 15626020
$ws.Range("I107").Value = 17858122
$ws.Range("K107").Value = 17858122
$ws.Range("M107").Value = -17856202

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2878.75
$ws.Range("I31").Value = 3151.7
$ws.Range("K31").Value = 3151.7
$ws.Range("M31").Value = -2856.7
$ws.Range("H34").Value = 2878.75
$ws.Range("I34").Value = 3151.7
$ws.Range("K34").Value = 3151.7
$ws.Range("M34").Value = -2949.7
$ws.Range("H58").Value = 4408.9355
$ws.Range("J58").Value = 8346.125
$ws.Range("L58").Value = 8346.125
$ws.Range("N58").Value = -8752.125
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
$ws.Range("H99").Value = 22008.4
$ws.Range("I99").Value = 22008.4
$ws.Range("K99").Value = 22008.4
$ws.Range("M99").Value = -20510.4
$ws.Range("H105").Value = 8437.3125
$ws.Range("I105").Value = 8928.357
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 8928.357
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -7181.357
$ws.Range("N105").Value = -8494
$ws.Range("H126").Value = 22008.4
$ws.Range("I126").Value = 22008.4
$ws.Range("K126").Value = 66025.20000000001
$ws.Range("M126").Value = -63555.20000000001
$ws.Range("H132").Value = 2894.8572
$ws.Range("I132").Value = 2894.8572
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8684.571599999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6154.571599999999
$ws.Range("N132").Value = $null
$ws.Range("H134").Value = 2305.0588
$ws.Range("I134").Value = 2323.25
$ws.Range("K134").Value = 6969.75
$ws.Range("M134").Value = -4434.75
$ws.Range("H136").Value = 4408.9355
$ws.Range("J136").Value = 8346.125
$ws.Range("L136").Value = 25038.375
$ws.Range("N136").Value = -30138.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 13102.083
$ws.Range("I141").Value = 9603.571
$ws.Range("K141").Value = 28810.713
$ws.Range("M141").Value = -23630.713

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 7219.6665
$ws.Range("I21").Value = 6247.75
$ws.Range("K21").Value = 6247.75
$ws.Range("M21").Value = -6074.75
$ws.Range("H30").Value = 7219.6665
$ws.Range("I30").Value = 6247.75
$ws.Range("K30").Value = 6247.75
$ws.Range("M30").Value = -6142.75
$ws.Range("H102").Value = 2178.1482
$ws.Range("I102").Value = 2444.2856
$ws.Range("J102").Value = 1246.6666
$ws.Range("K102").Value = 2444.2856
$ws.Range("L102").Value = 1246.6666
$ws.Range("M102").Value = -822.2856000000002
$ws.Range("N102").Value = -4490.6666
$ws.Range("H107").Value = 290.05884
$ws.Range("I107").Value = 266.86667
$ws.Range("J107").Value = 464
$ws.Range("K107").Value = 266.86667
$ws.Range("L107").Value = 464
$ws.Range("M107").Value = 1653.13333
$ws.Range("N107").Value = -4304
$ws.Range("H126").Value = 2722.6428
$ws.Range("I126").Value = 2427.375
$ws.Range("J126").Value = 3116.3333
$ws.Range("K126").Value = 7282.125
$ws.Range("L126").Value = 9348.999899999999
$ws.Range("M126").Value = -4812.125
$ws.Range("N126").Value = -14288.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 450
$ws.Range("I16").Value = 400
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 400
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -230
$ws.Range("N16").Value = -840
$ws.Range("H46").Value = 4353.522
$ws.Range("J46").Value = 4867.0527
$ws.Range("L46").Value = 4867.0527
$ws.Range("N46").Value = -5243.0527
$ws.Range("H55").Value = 1109.1777
$ws.Range("I55").Value = 952.15
$ws.Range("K55").Value = 952.15
$ws.Range("M55").Value = -779.15
$ws.Range("H132").Value = 2745.4443
$ws.Range("I132").Value = 2054.85
$ws.Range("J132").Value = 4718.5713
$ws.Range("K132").Value = 6164.549999999999
$ws.Range("L132").Value = 14155.7139
$ws.Range("M132").Value = -3634.549999999999
$ws.Range("N132").Value = -19215.7139

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5737.517
$ws.Range("J62").Value = 5626.7334
$ws.Range("L62").Value = 5626.7334
$ws.Range("N62").Value = -6874.7334
$ws.Range("H65").Value = 5737.517
$ws.Range("J65").Value = 5626.7334
$ws.Range("L65").Value = 28133.667
$ws.Range("N65").Value = -34373.667
$ws.Range("H81").Value = 3875.875
$ws.Range("I81").Value = 2636
$ws.Range("J81").Value = 5942.3335
$ws.Range("K81").Value = 5272
$ws.Range("L81").Value = 11884.667
$ws.Range("M81").Value = -4211
$ws.Range("N81").Value = -14006.667
$ws.Range("H84").Value = 3875.875
$ws.Range("I84").Value = 2636
$ws.Range("J84").Value = 5942.3335
$ws.Range("K84").Value = 26360
$ws.Range("L84").Value = 59423.335
$ws.Range("M84").Value = -21056
$ws.Range("N84").Value = -70031.33499999999
$ws.Range("H107").Value = 1960.3334
$ws.Range("I107").Value = 1234.8182
$ws.Range("J107").Value = 3411.3635
$ws.Range("K107").Value = 3704.4546
$ws.Range("L107").Value = 10234.0905
$ws.Range("M107").Value = -1784.4546
$ws.Range("N107").Value = -14074.0905
$ws.Range("H113").Value = 27778194
$ws.Range("J113").Value = 55555724
$ws.Range("L113").Value = 166667172
$ws.Range("N113").Value = -166671512
$ws.Range("H132").Value = 3350.125
$ws.Range("I132").Value = 2633.8125
$ws.Range("J132").Value = 4066.4375
$ws.Range("K132").Value = 7901.4375
$ws.Range("L132").Value = 12199.3125
$ws.Range("M132").Value = -5371.4375
$ws.Range("N132").Value = -17259.3125
$ws.Range("H136").Value = 1514.4546
$ws.Range("I136").Value = 1348.1034
$ws.Range("J136").Value = 2720.5
$ws.Range("K136").Value = 4044.3102
$ws.Range("L136").Value = 8161.5
$ws.Range("M136").Value = -1494.3102
$ws.Range("N136").Value = -13261.5
